# Reformat the JSON-ish payload (currently in A2, using Python-literal
# single-quote style) into pretty-printed, double-quoted JSON and move it
# up into A1 (replacing the old numeric placeholder / bold+bordered style
# that lived there), so the sheet ends up with a single used cell A1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
[
    {
        "title": "To connect your app to an existing database, you check the resource summary, as shown in the image below. You can see multiple resources. Which resources in the summary are database aliases?",
        "ques_type": 2,
        "options": [
            "Gaming-stars and gaming-stars-bot",
            "gaming-stars-bot-cloudantNoSQLDB and gaming-stars-bot-cloudantNoSQLDB",
            "Cloudant-pb and Cloudant-xl",
            "Watson Assistant-gv and Watson Discovery-o9"
        ],
        "score": "gaming-stars-bot-cloudantNoSQLDB and gaming-stars-bot-cloudantNoSQLDB"
    },
    {
        "title": "Your company is using IBM Watson Assistant to manage customer support questions for one of its products. Management would like to add the same solution to five other products. However, the questions will be mostly the same.  How would you scale IBM Watson Assistant?",
        "ques_type": 2,
        "options": [
            "Create a new instance for each product.",
            "Use the same instance and create a new assistant for each product.",
            "Add a context variable so the existing assistant knows which product is currently discussed.",
            "Create a new dialog node branch for each product."
        ],
        "score": "Add a context variable so the existing assistant knows which product is currently discussed."
    },
    {
        "title": "You are leading a meeting about how your client can utilize IBM Watson in their company. The company is an online retailer selling premium accessories. They have asked you to create several modules for them, and they request that you let them know which of their desired modules are impossible to make using IBM Watson.  Which of these modules are outside the scope of IBM Watson functionality?",
        "ques_type": 2,
        "options": [
            "A customer support chatbot.",
            "A cognitive search engine to search for products.",
            "A search engine that finds products based on similar images.",
            "A machine-learning model that recommends products."
        ],
        "score": "A search engine that finds products based on similar images."
    },
    {
        "title": "You are building an IBM Watson Assistant chatbot for your company. The chatbot has to be able to give a short description of each of the company's 10 employees. What is the fastest way to create this functionality?",
        "ques_type": 2,
        "options": [
            "Create an entity @employee and a dialog node with a story for each employee.",
            "Connect to a database with employees and retrieve stories based on a wildcard search.",
            "Connect to your company\u2019s website and scrape all empty data.",
            "Connect to LinkedIn API and search for the company\u2019s employees and their descriptions there."
        ],
        "score": "Create an entity @employee and a dialog node with a story for each employee."
    }
]
'@

# A2 (old shared-string cell) is emptied - its content now lives in A1.
$ws.Range("A2").ClearContents()
$ws.Range("A2").ClearFormats()

# A1 (old bold/bordered numeric placeholder cell) loses its formatting and
# gets the new text instead.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = "questions = " + $text

# Setting a value containing embedded newlines auto-expands the row height;
# AutoFit puts the row back to the (non-custom) default height.
$ws.Rows.Item(1).AutoFit()
